# Generate Report for Handoff
# - Flip Status from "In Translation" to "Ready for handoff" on all three
#   sheets (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to the new handoff run time.
# - Widen the (now longer) status columns to fit "Ready for handoff".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps for the new handoff generation ---
$wsOverview.Range("G2").Value = "2016-08-21 16:47:20"
$wsDeDe.Range("H2").Value = "2016-08-21 16:47:20"
$wsZhCn.Range("H2").Value = "2016-08-21 16:47:16"

# --- Widen the Status columns now that the text is longer ---
$wsOverview.Range("E1:F1").ColumnWidth = 16.38265482584637
$wsZhCn.Range("C1").ColumnWidth = 16.38265482584637
$wsDeDe.Range("C1").ColumnWidth = 16.38265482584637
